$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.305.97'
$ws.Range("E2").Value = '  -0.12%  '

$ws.Range("D3").Value = '3.135.08'
$ws.Range("E3").Value = '  -1.39%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '''571.03'
$ws.Range("E5").Value = '  -0.09%  '

$ws.Range("D6").Value = '''163.69'
$ws.Range("E6").Value = '  -2.98%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").Value = '''0.573'
$ws.Range("E8").Value = '  -5.84%  '

$ws.Range("D9").Value = '3.147.89'
$ws.Range("E9").Value = '  -1.21%  '

$ws.Range("E10").Value = '  -3.26%  '

$ws.Range("E11").Value = '  -2.62%  '

$ws.Range("E12").Value = '  -1.28%  '

$ws.Range("D13").Value = '3.686.76'
$ws.Range("E13").Value = '  -1.19%  '

$ws.Range("E14").Value = '  -1.75%  '

$ws.Range("D15").Value = '64.265.48'
$ws.Range("E15").Value = '  -0.27%  '

$ws.Range("D16").Value = '''24.93'
$ws.Range("E16").Value = '  -1.59%  '

$ws.Range("D17").Value = '3.143.79'
$ws.Range("E17").Value = '  -1.10%  '

$ws.Range("E18").Value = '  -2.69%  '

$ws.Range("D19").Value = '''405.42'
$ws.Range("E19").Value = '  -3.13%  '

$ws.Range("E20").Value = '  -2.27%  '

$ws.Range("E21").Value = '  -3.51%  '

$ws.Range("E22").Value = '  -0.61%  '

$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  +0.08%  '

$ws.Range("D24").Value = '''68.94'
$ws.Range("E24").Value = '  -2.06%  '

$ws.Range("D25").Value = '''0.483'
$ws.Range("E25").Value = '  -1.11%  '

$ws.Range("E26").Value = '  -5.31%  '

$ws.Range("E27").Value = '  -3.85%  '

$ws.Range("D28").Value = '''8.86'
$ws.Range("E28").Value = '  +1.31%  '

$ws.Range("D29").Value = '''0.996'
$ws.Range("E29").Value = '  -0.10%  '

$ws.Range("E30").Value = '  +0.08%  '

$ws.Range("E31").Value = '  -2.21%  '

$ws.Range("D32").Value = '''21.24'
$ws.Range("E32").Value = '  -2.34%  '

$ws.Range("D33").Value = '''162.14'
$ws.Range("E33").Value = '  +3.40%  '

$ws.Range("D34").Value = '''4.86'
$ws.Range("E34").Value = '  -3.81%  '

$ws.Range("D35").Value = '''6.27'
$ws.Range("E35").Value = '  -1.22%  '

$ws.Range("D36").Value = '''1.12'
$ws.Range("E36").Value = '  -0.84%  '

$ws.Range("E37").Value = '  -0.85%  '

$ws.Range("E38").Value = '  -1.50%  '

$ws.Range("D39").Value = '2.635.03'
$ws.Range("E39").Value = '  -2.69%  '

$ws.Range("D40").Value = '''23.58'
$ws.Range("E40").Value = '  -2.82%  '

$ws.Range("D41").Value = '''4.08'
$ws.Range("E41").Value = '  -3.30%  '

$ws.Range("D42").Value = '''38.32'
$ws.Range("E42").Value = '  -2.29%  '

$ws.Range("D43").Value = '''0.690'
$ws.Range("E43").Value = '  -3.80%  '

$ws.Range("D44").Value = '''0.0611'
$ws.Range("E44").Value = '  -1.98%  '

$ws.Range("D45").Value = '''5.40'
$ws.Range("E45").Value = '  -3.05%  '

$ws.Range("D46").Value = '''287.87'
$ws.Range("E46").Value = '  -1.68%  '

$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").Value = '''0.0254'
$ws.Range("E47").Value = '  -4.11%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").Value = '''21.12'
$ws.Range("E48").Value = '  -1.43%  '

$ws.Range("D49").Value = '''0.997'
$ws.Range("E49").Value = '  -0.09%  '

$ws.Range("D50").Value = '''0.0975'
$ws.Range("E50").Value = '  -1.62%  '

$ws.Range("E51").Value = '  +0.57%  '
